$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 25: add H25 and I25 (availability hours for the 4th day, afternoon columns)
$ws.Range("H25").Value = 0.083333333333333329
$ws.Range("H25").NumberFormat = $ws.Range("F25").NumberFormat

$ws.Range("I25").Value = 0.027777777777777776
$ws.Range("I25").NumberFormat = $ws.Range("F25").NumberFormat

# Row 29: new task "modifierdisponibilites"
$ws.Range("A29").Value = "modifierdisponibilites"

$ws.Range("H29").Value = 0.041666666666666664
$ws.Range("H29").NumberFormat = $ws.Range("F25").NumberFormat

$ws.Range("I29").Value = 0.125
$ws.Range("I29").NumberFormat = $ws.Range("F25").NumberFormat

# Row 30: new task "Commentaires"
$ws.Range("A30").Value = "Commentaires"

$ws.Range("H30").Value = 0.03125
$ws.Range("H30").NumberFormat = $ws.Range("F25").NumberFormat

# Update the selected cell to match the new active selection
$ws.Range("I25").Select()
